$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Row 2 ----
$ws.Range("E2").Value = 23.98000000000031
$ws.Range("H2").Value = [double]"1.22170346588738e-16"
$ws.Range("K2").Value = 51.1346574274139
$ws.Range("L2").Value = "[45.58196413732337, 56.68735071750443]"
$ws.Range("O2").Value = 1.679289766783733
$ws.Range("P2").Value = "[1.566079220708425, 1.7925003128590413]"
$ws.Range("S2").Value = 57.41933532847631
$ws.Range("T2").Value = "[53.85399752241417, 60.98467313453846]"
$ws.Range("W2").Value = 17.57093093093116
$ws.Range("X2").Value = 17.13885885885908
$ws.Range("Y2").Value = 18.00300300300324

# ---- Row 3 ----
$ws.Range("E3").Value = 23.66000000000026
$ws.Range("H3").Value = [double]"1.22170346588738e-16"
# p_reject (I3) becomes blank/not-significant, same empty-text shape as I2
$i3Style = $ws.Range("I3").Style
$ws.Range("I3").Value = "'"
$ws.Range("I3").Style = $i3Style
$ws.Range("K3").Value = 52.97807066987878
$ws.Range("L3").Value = "[43.55707608422088, 62.39906525553667]"
$ws.Range("O3").Value = 0.9119738433844251
$ws.Range("P3").Value = "[0.723289599925578, 1.1006580868432723]"
$ws.Range("S3").Value = 57.76995394345558
$ws.Range("T3").Value = "[52.888522590455146, 62.651385296456006]"
$ws.Range("W3").Value = 20.22586586586608
$ws.Range("X3").Value = 19.51535535535556
$ws.Range("Y3").Value = 20.9363763763766
